{"js": "// The campaign-dates line currently reads (Spanish):\n//   \"Fechas de la campa\u00f1a para Constelaci\u00f3n de H\u00e9rcules 2022: 13-22 de junio, 12-21 de julio, 10-19 de agosto\"\n// and needs to become:\n//   \"2022 Fechas de la campa\u00f1a para Constelaci\u00f3n de H\u00e9rcules: 13-22 de junio, 12-21 de julio, 10-19 de agosto\"\n// (moving the leading \"2022\" year to the front of the sentence). This exact\n// sentence appears 4 times in the document body, each time as the sole text\n// of its own run/paragraph, so we find every occurrence and swap it in place.\n\nconst oldText =\n  \"Fechas de la campa\u00f1a para Constelaci\u00f3n de H\u00e9rcules 2022: 13-22 de junio, 12-21 de julio, 10-19 de agosto\";\nconst newText =\n  \"2022 Fechas de la campa\u00f1a para Constelaci\u00f3n de H\u00e9rcules: 13-22 de junio, 12-21 de julio, 10-19 de agosto\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The campaign-dates line currently reads (Spanish):\n#   \"Fechas de la campa\u00f1a para Constelaci\u00f3n de H\u00e9rcules 2022: 13-22 de junio, 12-21 de julio, 10-19 de agosto\"\n# and needs to become:\n#   \"2022 Fechas de la campa\u00f1a para Constelaci\u00f3n de H\u00e9rcules: 13-22 de junio, 12-21 de julio, 10-19 de agosto\"\n# (moving the leading \"2022\" year to the front of the sentence). This exact\n# sentence appears 4 times throughout the document body, so use Find/Replace\n# across the whole document (wdReplaceAll) to update every occurrence.\n\n$d = $word.ActiveDocument\n\n$oldText = \"Fechas de la campa\u00f1a para Constelaci\u00f3n de H\u00e9rcules 2022: 13-22 de junio, 12-21 de julio, 10-19 de agosto\"\n$newText = \"2022 Fechas de la campa\u00f1a para Constelaci\u00f3n de H\u00e9rcules: 13-22 de junio, 12-21 de julio, 10-19 de agosto\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n\n# Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n#         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n#   Wrap = 1 (wdFindContinue), Replace = 2 (wdReplaceAll)\n$find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n"}
